$wb = $excel.ActiveWorkbook

# Sheet "展览" updates (想去人数 column F)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 49
$ws1.Range("F7").Value = 578
$ws1.Range("F8").Value = 62
$ws1.Range("F9").Value = 8396
$ws1.Range("F10").Value = 782
$ws1.Range("F12").Value = 1126
$ws1.Range("F13").Value = 893
$ws1.Range("F17").Value = 158
$ws1.Range("F19").Value = 222
$ws1.Range("F20").Value = 935

# Sheet "全部类型" updates (想去人数 column F)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 49
$ws4.Range("F9").Value = 578
$ws4.Range("F10").Value = 62
$ws4.Range("F11").Value = 8396
$ws4.Range("F12").Value = 782
$ws4.Range("F14").Value = 1126
$ws4.Range("F15").Value = 893
$ws4.Range("F19").Value = 158
$ws4.Range("F21").Value = 222
$ws4.Range("F22").Value = 935
